$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 values (metricas_recorrencia_anual)
$ws.Range("C8").Value = 1181
$ws.Range("E8").Value = 988
$ws.Range("G8").Value = 83.65791701947502
$ws.Range("H8").Value = 16.34208298052498
